# New weekly price record for "Ajo" (garlic) at Terminal Hortofrutícola Agro
# Chillán: insert a new row above the current row 17 (pushing the existing
# rows 17-126 down to 18-127) and populate it with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(17).Insert()

$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44473
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112003
$ws.Range("G17").Value = "Ajo"
$ws.Range("H17").Value = "Chino"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 16000
$ws.Range("L17").Value = 17000
$ws.Range("M17").Value = 16500
$ws.Range("N17").Value = "$/caja 10 kilos"
$ws.Range("O17").Value = "China"
$ws.Range("P17").Value = 1650
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = "Hortaliza"
